$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a brand-new match row at row 5 (this shifts the former rows 5 and
#    6 down to rows 6 and 7 respectively).
# ---------------------------------------------------------------------------
$ws.Rows(5).Insert()

# Fill in the data for the newly inserted row (Netherlands - Tweede Divisie
# match between Feyenoord U21 and Jong Sparta Rotterdam). Odds are not yet
# available for this fixture, so the numeric columns stay blank.
$ws.Range("A5").Value = "nDjP3tNH"
# Leading apostrophe forces this date-shaped token to stay literal text
# (matching the other rows' Date column) instead of being parsed into a
# date serial number.
$ws.Range("B5").Value = "'03/06/2025"
$ws.Range("C5").Value = "14:30"
$ws.Range("D5").Value = "NETHERLANDS - TWEEDE DIVISIE"
$ws.Range("E5").Value = "Feyenoord U21"
$ws.Range("F5").Value = "Jong Sparta Rotterdam"

# ---------------------------------------------------------------------------
# 2) Row 2 (U. Espanola vs Limache) - updated odds
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 2.15
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 15
$ws.Range("P2").Value = 1.29
$ws.Range("Q2").Value = 3.5
$ws.Range("U2").Value = 13
$ws.Range("Z2").Value = 15

# ---------------------------------------------------------------------------
# 3) Row 3 (Nueve de Octubre vs Leones del Norte) - updated odds
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = 2.5
$ws.Range("H3").Value = 2.75
$ws.Range("I3").Value = 3.05
$ws.Range("L3").Value = 1.52
$ws.Range("M3").Value = 2.2
$ws.Range("N3").Value = 2.47
$ws.Range("O3").Value = 1.42
$ws.Range("P3").Value = 1.55
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 5.8
$ws.Range("U3").Value = 10.75
$ws.Range("V3").Value = 10.5
$ws.Range("W3").Value = 28
$ws.Range("X3").Value = 28
$ws.Range("Y3").Value = 50
$ws.Range("Z3").Value = 5.9
$ws.Range("AA3").Value = 5.6
$ws.Range("AB3").Value = 17.5
$ws.Range("AC3").Value = 120
$ws.Range("AE3").Value = 7.1
$ws.Range("AF3").Value = 14.5
$ws.Range("AG3").Value = 11.5
$ws.Range("AH3").Value = 40
$ws.Range("AI3").Value = 32
$ws.Range("AJ3").Value = 50

# ---------------------------------------------------------------------------
# 4) Row 4 (Imbabura vs Chacaritas) - updated odds (AD4 now has a value)
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 1.45
$ws.Range("H4").Value = 3.85
$ws.Range("I4").Value = 7.2
$ws.Range("L4").Value = 1.34
$ws.Range("M4").Value = 2.75
$ws.Range("N4").Value = 1.98
$ws.Range("O4").Value = 1.65
$ws.Range("P4").Value = 1.42
$ws.Range("Q4").Value = 2.47
$ws.Range("T4").Value = 5.3
$ws.Range("U4").Value = 5.8
$ws.Range("V4").Value = 8.5
$ws.Range("W4").Value = 9.25
$ws.Range("X4").Value = 13.5
$ws.Range("Y4").Value = 37
$ws.Range("Z4").Value = 8.5
$ws.Range("AA4").Value = 7.8
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 15.5
$ws.Range("AF4").Value = 45
$ws.Range("AG4").Value = 24
$ws.Range("AH4").Value = 200
$ws.Range("AI4").Value = 100
$ws.Range("AJ4").Value = 100

# ---------------------------------------------------------------------------
# 5) Row 6 (Tacuarembo vs Cerrito, formerly row 5) is unchanged by the shift.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 6) Row 7 (Rampla Juniors vs Oriental, formerly row 6) - updated odds
# ---------------------------------------------------------------------------
$ws.Range("G7").Value = 2.6
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 2.47
$ws.Range("L7").Value = 1.32
$ws.Range("M7").Value = 2.82
$ws.Range("N7").Value = 1.93
$ws.Range("O7").Value = 1.7
$ws.Range("P7").Value = 1.39
$ws.Range("Q7").Value = 2.57
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.83
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 12.5
$ws.Range("W7").Value = 28
$ws.Range("Y7").Value = 35
$ws.Range("Z7").Value = 9.25
$ws.Range("AA7").Value = 6.4
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 75
$ws.Range("AD7").Value = 600
$ws.Range("AE7").Value = 7.8
$ws.Range("AG7").Value = 9.75
$ws.Range("AH7").Value = 26
$ws.Range("AI7").Value = 21
$ws.Range("AJ7").Value = 32
